# Apply edit: add two new quarterly columns (D, E) to the MSFT financials sheet,
# shifting existing data (old D:K) right to (F:M), and correct a handful of
# previously-mis-entered figures for the 2017-06-30 quarter (old column H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert two new blank columns at D:E. This shifts old D:K -> F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) Copy number/date formatting from column F (first untouched data column)
#    into the two new columns so the new cells pick up the same styles
#    (date format for row type 7/38/80, number format elsewhere).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill in the two new columns (D = newest quarter, E = prior quarter)
#    for every row that carries data.
$newQuarterData = @{
    7 = @(43465, 43373)
    8 = @(32471000, 29084000)
    9 = @(12423000, 9905000)
    10 = @(20048000, 19179000)
    12 = @(4070000, 3977000)
    13 = @(0, 0)
    14 = @(7000, 0)
    15 = @(0, 0)
    17 = @(22220000, 19129000)
    18 = @(10251000, 9955000)
    20 = @(806000, 940000)
    21 = @(14052000, 13732000)
    22 = @(672000, 674000)
    23 = @(10385000, 10221000)
    24 = @(1965000, 1397000)
    25 = @(0, 0)
    26 = @(8420000, 8824000)
    27 = @(8420000, 8824000)
    28 = @(0, 0)
    29 = @("NA", "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-806000, -940000)
    33 = @(8420000, 8824000)
    34 = @(0, 0)
    35 = @(8420000, 8824000)
    38 = @(43465, 43373)
    41 = @(6638000, 15137000)
    42 = @(121013000, 120634000)
    43 = @(19680000, 17390000)
    44 = @(1961000, 3614000)
    45 = @(7582000, 7420000)
    46 = @(156874000, 164195000)
    47 = @(2274000, 2034000)
    48 = @(39523000, 38164000)
    49 = @(50059000, 43434000)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(10129000, 9792000)
    53 = @(0, 0)
    54 = @(258859000, 257619000)
    57 = @(7563000, 8511000)
    58 = @(3761000, 6700000)
    59 = @(38994000, 41066000)
    60 = @(50318000, 56277000)
    61 = @(74618000, 74204000)
    62 = @(41795000, 41171000)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(166731000, 171652000)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(16585000, 17279000)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(92128000, 85967000)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(8420000, 8824000)
    83 = @(2995000, 2837000)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(8900000, 13657000)
    91 = @(-3707000, -3602000)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-4200000, -2953000)
    96 = @(-3544000, -3220000)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-13216000, -7384000)
    101 = @(17000, -129000)
    102 = @(-8499000, 3191000)
}

foreach ($row in $newQuarterData.Keys) {
    $vals = $newQuarterData[$row]
    $ws.Cells.Item([int]$row, 4).Value2 = $vals[0]
    $ws.Cells.Item([int]$row, 5).Value2 = $vals[1]
}

# 4) A few rows also received corrected historical figures (2017-06-30
#    quarter, now column J) as part of this update, not just a shift of the
#    old value.
$corrections = @{
    20 = @{"J"=1587000}
    21 = @{"J"=11611000}
    22 = @{"J"=1310000}
    32 = @{"J"=-1587000}
    48 = @{"F"=36146000; "G"=34788000; "H"=33053000; "I"=31653000; "J"=30289000}
    49 = @{"F"=43736000; "G"=44126000; "H"=44389000; "I"=44987000; "J"=45228000}
    59 = @{"J"=38121000}
    91 = @{"J"=-2283000}
}

foreach ($row in $corrections.Keys) {
    $colmap = $corrections[$row]
    foreach ($col in $colmap.Keys) {
        $ws.Range("$col$row").Value2 = $colmap[$col]
    }
}

Write-Host "edit applied"
